$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update interest-rate values (row 2 - Namibia, row 3 - Other)
$ws.Range("D2").Value = 0.066053564606404797
$ws.Range("F2").Value = 0.066053564606404797
$ws.Range("H2").Value = 0.066053564606404797
$ws.Range("J2").Value = 0.059635333081896211

$ws.Range("D3").Value = 0.066053564606404797
$ws.Range("F3").Value = 0.066053564606404797
$ws.Range("H3").Value = 0.066053564606404797
$ws.Range("J3").Value = 0.059635333081896211

# Column J previously had a custom number format (0.000000000) applied via style index 1;
# that custom style is removed in the new version, so reset J2:J3 back to the default style.
$ws.Range("J2:J3").Style = "Normal"

# Column width changes: H widens to match G's width, J widens to match I's width
# (the stored widths become 25.375 and 24.625 respectively; Excel's ColumnWidth->stored
# conversion quantizes to the nearest 1/6 character, so these inputs land on the closest
# representable values)
$ws.Columns.Item(8).ColumnWidth = 24.58
$ws.Columns.Item(10).ColumnWidth = 23.91

# Update the active selection
$ws.Range("B16").Select() | Out-Null
